# Apply data corrections to the "dSF" column (F) for specific rows,
# per commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -5
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = -4
$ws.Range("F17").Value = 2
